$wb = $excel.ActiveWorkbook

# Values to update: column F (想去人数) for rows 2-5
$updates = @{
    2 = 10145
    3 = 226
    4 = 52
    5 = 623
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
